$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last column (BA), which held the "Run 50" header/data and the old
# "Mean" column. This shifts nothing else, it just drops the last column and
# shrinks the used range from A1:BA14 to A1:AZ14.
$ws.Columns("BA").Delete()

# Rename the "Gen" header (A1) to "MaxFES"
$ws.Range("A1").Value = "MaxFES"

# Replace the Gen values in column A with the new MaxFES fractional values
$maxFesValues = @(0, 0.001, 0.01, 0.1, 0.2, 0.3, 0.4, 0.5, 0.6, 0.7, 0.8, 0.9, 1)
for ($i = 0; $i -lt $maxFesValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $maxFesValues[$i]
}

# AZ1 still holds the label for the last data column; it should now read "Mean"
# (previously it was "Run 50", now Run 50 has been dropped and Mean moved here).
$ws.Range("AZ1").Value = "Mean"

# Recompute the Mean column (now in AZ) using only the 50 remaining runs
# (columns B:AY), since the Run 50 column has been removed.
$meanValues = @(
    367490228.0594656,
    260635150.1331378,
    103098732.5721272,
    14527053.80068066,
    7148699.11142516,
    4029228.05043914,
    2342353.62845974,
    1186599.34047946,
    679454.08566992,
    389657.27424091,
    277410.87852327,
    148714.91343774,
    83314.46526154
)
for ($i = 0; $i -lt $meanValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 52).Value = $meanValues[$i]
}
